$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D: "Locator Type" header and "CSS" values for the 5 data rows
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Locator Type"

$ws.Range("D2:D6").Value = "CSS"

# Update the selection to match the new active range
$ws.Range("D2:D6").Select()
